$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C, rows 2 through 114 hold the "Förändrad" (changed) date.
# All of them currently store the serial date value 45171 and must be
# bumped by one day to 45172.
$ws.Range("C2:C114").Value = 45172
